$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.368.86"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.602.60"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.03"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.57"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "2.626.27"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.328"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.131"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "3.060.27"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "58.316.49"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.58"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "2.591.34"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.61"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.36"
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.47"
$ws.Range("E24").Value = "  +4.28%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.405"
$ws.Range("E26").Value = "  -4.01%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.714.87"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.07"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "0.0₃0756"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.28"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.87"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.69"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.07"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.850"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.18"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.58"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.26"
$ws.Range("E43").Value = "  +3.77%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.99"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0525"
$ws.Range("E49").Value = "  -1.12%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.21"
$ws.Range("E50").Value = "  +5.73%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.74"
$ws.Range("E51").Value = "  +2.60%  "
